$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Remove the old merged header cells before reshaping the table ---
$ws.Range("E1:G1").UnMerge()
$ws.Range("H1:J1").UnMerge()
$ws.Range("K1:L1").UnMerge()

# --- Row 1: top-level headers ---
$ws.Range("A1").Value = "Start date"
$ws.Range("B1").Value = "End date"
$ws.Range("C1").Value = "Litter bag 1 (g)"
$ws.Range("D1").Value = "Litter bag 2 (g)"
$ws.Range("E1").Value = "Litter bag 3 (g)"
$ws.Range("J1").Value = "Litter bag 4 (g)"
$ws.Range("O1").Value = "Litter bag 5/10 (g)"
$ws.Range("R1").Value = "Litter bag 7 (g)"

# --- Row 2: species sub-headers ---
$ws.Range("C2").Value = "Oak"
$ws.Range("D2").Value = "Oak"
$ws.Range("E2").Value = "Oak "
$ws.Range("F2").Value = "Hazel"
$ws.Range("G2").Value = "Beech"
$ws.Range("H2").Value = "Ash"
$ws.Range("I2").Value = "Willow"
$ws.Range("J2").Value = "Oak"
$ws.Range("K2").Value = "Hazel"
$ws.Range("L2").Value = "Ash"
$ws.Range("M2").Value = "Hawthorn"
$ws.Range("N2").Value = "Willow"
$ws.Range("O2").Value = "Oak"
$ws.Range("P2").Value = "Hazel"
$ws.Range("Q2").Value = "Pine"
$ws.Range("R2").Value = "Oak"
$ws.Range("S2").Value = "Hazel"
$ws.Range("T2").Value = "Other"

# --- Row 3: first period of data (replaces old narrower row) ---
$ws.Range("A3").Value = 42233
$ws.Range("B3").Value = 42255
$ws.Range("C3").Value = 1.74
$ws.Range("D3").Value = 1.1000000000000001
$ws.Range("E3").Value = 0.72
$ws.Range("F3").Value = 0.73
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0.03
$ws.Range("J3").Value = 0.54
$ws.Range("K3").Value = 0.7
$ws.Range("L3").Value = 0.54
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 0
$ws.Range("O3").Value = 1.05
$ws.Range("P3").Value = 0.2
$ws.Range("Q3").Value = 0
$ws.Range("R3").Value = 2.08
$ws.Range("S3").Value = 0
$ws.Range("T3").Value = 0.1

# --- Row 4 ---
$ws.Range("A4").Value = 42255
$ws.Range("B4").Value = 42275
$ws.Range("C4").Value = 1.6
$ws.Range("D4").Value = 1.57
$ws.Range("E4").Value = 1.1100000000000001
$ws.Range("F4").Value = 1.18
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0.69
$ws.Range("K4").Value = 0.61
$ws.Range("L4").Value = 0.35
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0.28999999999999998
$ws.Range("O4").Value = 1.97
$ws.Range("P4").Value = 0.18
$ws.Range("Q4").Value = 0.5
$ws.Range("R4").Value = 0.89
$ws.Range("S4").Value = 0
$ws.Range("T4").Value = 0

# --- Row 5 ---
$ws.Range("A5").Value = 42275
$ws.Range("B5").Value = 42292
$ws.Range("C5").Value = 5.25
$ws.Range("D5").Value = 5.12
$ws.Range("E5").Value = 4.84
$ws.Range("F5").Value = 2.93
$ws.Range("G5").Value = 0.44
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0.44
$ws.Range("J5").Value = 1.77
$ws.Range("K5").Value = 3.2
$ws.Range("L5").Value = 5.36
$ws.Range("M5").Value = 0.27
$ws.Range("N5").Value = 0.26
$ws.Range("O5").Value = 6.14
$ws.Range("P5").Value = 0.73
$ws.Range("Q5").Value = 0.26
$ws.Range("R5").Value = 4.55
$ws.Range("S5").Value = 0.3
$ws.Range("T5").Value = 0

# --- Row 6 ---
$ws.Range("A6").Value = 42292
$ws.Range("B6").Value = 42306
$ws.Range("C6").Value = 11.73
$ws.Range("D6").Value = 12.47
$ws.Range("E6").Value = 12.14
$ws.Range("F6").Value = 2.44
$ws.Range("G6").Value = 0.4
$ws.Range("H6").Value = 1.02
$ws.Range("I6").Value = 0.4
$ws.Range("J6").Value = 8.35
$ws.Range("K6").Value = 11.34
$ws.Range("L6").Value = 4.7300000000000004
$ws.Range("M6").Value = 2.23
$ws.Range("N6").Value = 0
$ws.Range("O6").Value = 10.67
$ws.Range("P6").Value = 1.84
$ws.Range("Q6").Value = 2.08
$ws.Range("R6").Value = 12.03
$ws.Range("S6").Value = 0
$ws.Range("T6").Value = 0

# --- Number formats ---
$ws.Range("A3:B6").NumberFormat = "mm/dd/yyyy"

# --- Re-merge the header groups across the new wider spans ---
$ws.Range("E1:I1").Merge()
$ws.Range("J1:N1").Merge()
$ws.Range("O1:Q1").Merge()
$ws.Range("R1:T1").Merge()

# Center the merged header row
$ws.Range("E1:T1").HorizontalAlignment = -4108

# --- Column widths ---
$ws.Columns.Item(1).ColumnWidth = 12
$ws.Columns.Item(2).ColumnWidth = 12.28515625

# --- Sheet view: zoom and selection ---
$ws.Application.ActiveWindow.Zoom = 90
$ws.Range("I20").Select()
